$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")
$spain = $wb.Worksheets.Item("Spain")
$czech = $wb.Worksheets.Item("Czech")

# --- Netherlands: copy from Italy (same column widths / row-height layout) ---
$italy.Copy($null, $spain)
$netherlands = $wb.Worksheets.Item(8)
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Value = "NGC-3144/T2199"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Range("B4").Select()

# --- Austria: copy from Czech, then drop the PR1D2-Unmonitored row ---
$czech.Copy($null, $netherlands)
$austria = $wb.Worksheets.Item(9)
$austria.Name = "Austria"
$austria.Range("B4").Value = "NGC-3817/T2306"
$austria.Range("B2").Value = "Austria Market"

$austria.Rows.Item(10).Delete()
$austria.Rows.Item(11).EntireRow.AutoFit()
$austria.Rows.Item(12).RowHeight = 13.8
$austria.Range("B17").Select()

# --- Denmark: copy from Czech ---
$czech.Copy($null, $austria)
$denmark = $wb.Worksheets.Item(10)
$denmark.Name = "Denmark"
$denmark.Range("B4").Value = "NGC-2913/T2798"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Select()

# --- Spain tab is no longer the active/selected tab ---
$spain.Range("G9").Select()

# --- Austria becomes the active sheet/tab ---
$austria.Activate()
